$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8513891100883484
$ws.Range("B1").Value = 1.361099243164062
$ws.Range("C1").Value = 2.240350484848022
$ws.Range("D1").Value = 2.302311897277832
$ws.Range("E1").Value = 1.793417811393738
